# Generate Report for Handoff
# Regenerates the handoff artifacts (new source GUID / new xliff hashes),
# refreshes the "Latest Handoff" file names + timestamps, and clears the
# "Latest Target / Latest Handback" columns since the freshly generated
# handoff has not been targeted or handed back yet.

$wb = $excel.ActiveWorkbook

$oldGuid = "67ca1668-05ae-4bbd-a9ae-6c043f570a87"
$newGuid = "13d586ca-986c-4be8-b7e3-cf69292a6666"
$newHash = "7e0eda1170f58dd895f860ebaaa0fb41183418e0"

$newFileName = "$newGuid.md"
$newHandoffDateZh = "2016-08-25 22:58:16"
$newHandoffDateDe = "2016-08-25 22:58:21"
$neverDate = "0001-01-01 00:00:00"
$newXlfZh = "$newGuid.$newHash.zh-cn.xlf"
$newXlfDe = "$newGuid.$newHash.de-de.xlf"

# NOTE: the handoff-file rename only touches the *display* text of each
# hyperlink; the underlying link target keeps pointing at the original
# (un-renamed) path in the source repo, so we re-create each hyperlink
# against the same old-guid address.
$oldFileName = "$oldGuid.md"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57d512d28145d03a5fd9e4934fcd37d128094ff2/e2e"
$overviewAddr = "$baseUrl/$oldFileName"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = "e2e\$newFileName"
$wsOverview.Range("G2").Value = $newHandoffDateDe

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewAddr, "", "", "e2e\$newFileName")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newHandoffDateZh
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("J2").Style = "Normal"
$wsZh.Range("K2").Value = $neverDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $overviewAddr, "", "", $newFileName)

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newHandoffDateDe
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("J2").Style = "Normal"
$wsDe.Range("K2").Value = $neverDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $overviewAddr, "", "", $newFileName)

# ---------------------------------------------------------------------
# Column widths for the now-empty "Latest Target File" / "Latest Handback
# File" columns (autofit down from the old 40-char fixed width).
# ---------------------------------------------------------------------
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
